$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: name + link text (no second URL column for this entry)
$ws.Range("A13").Value = "THE HUMAN PROTEIN ATLAS"
$ws.Range("B13").Value = "https://www.proteinatlas.org/"

# Move selection to A14, matching the post-edit state
$ws.Range("A14").Select()
